$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price/Volume columns so numeric-looking values
# (e.g. "0.694", "43.73") are stored as text, matching the source data
# which uses inline strings throughout (prices use "." as both a
# thousands separator and decimal point, so they must stay text).
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '36.310.70'
$ws.Range("E2").Value = '  +2.65%  '

# Row 3
$ws.Range("D3").Value = '1.916.10'
$ws.Range("E3").Value = '  +1.38%  '

# Row 4
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
$ws.Range("D5").Value = '249.08'
$ws.Range("E5").Value = '  +1.03%  '

# Row 6
$ws.Range("D6").Value = '0.694'
$ws.Range("E6").Value = '  +0.34%  '

# Row 7
$ws.Range("E7").Value = '  +0.08%  '

# Row 8
$ws.Range("D8").Value = '43.73'
$ws.Range("E8").Value = '  +1.08%  '

# Row 9
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").Value = '0.364'
$ws.Range("E9").Value = '  +2.88%  '

# Row 10
$ws.Range("B10").Value = 'OKB'
$ws.Range("C10").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D10").Value = '57.37'
$ws.Range("E10").Value = '  +6.50%  '

# Row 11
$ws.Range("D11").Value = '0.0763'
$ws.Range("E11").Value = '  +2.68%  '

# Row 12
$ws.Range("D12").Value = '0.0995'
$ws.Range("E12").Value = '  +2.51%  '

# Row 13
$ws.Range("E13").Value = '  +9.76%  '

# Row 14
$ws.Range("D14").Value = '0.796'
$ws.Range("E14").Value = '  +4.85%  '

# Row 15
$ws.Range("D15").Value = '2.201.79'
$ws.Range("E15").Value = '  +1.65%  '

# Row 16
$ws.Range("D16").Value = '5.11'
$ws.Range("E16").Value = '  +4.42%  '

# Row 17
$ws.Range("D17").Value = '1.913.32'
$ws.Range("E17").Value = '  +2.21%  '

# Row 18
$ws.Range("D18").Value = '36.376.50'
$ws.Range("E18").Value = '  +2.67%  '

# Row 19
$ws.Range("D19").Value = '74.36'
$ws.Range("E19").Value = '  +1.56%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0844'
$ws.Range("E20").Value = '  +2.35%  '

# Row 21
$ws.Range("D21").Value = '252.77'
$ws.Range("E21").Value = '  +3.16%  '

# Row 22
$ws.Range("D22").Value = '13.19'
$ws.Range("E22").Value = '  +3.00%  '

# Row 23
$ws.Range("D23").Value = '5.20'
$ws.Range("E23").Value = '  +4.93%  '

# Row 24
$ws.Range("D24").Value = '2.69'
$ws.Range("E24").Value = '  +0.37%  '

# Row 25
$ws.Range("E25").Value = '  -0.03%  '

# Row 26
$ws.Range("D26").Value = '2.26'
$ws.Range("E26").Value = '  +5.36%  '

# Row 27
$ws.Range("D27").Value = '167.58'
$ws.Range("E27").Value = '  +0.83%  '

# Row 28
$ws.Range("D28").Value = '8.80'
$ws.Range("E28").Value = '  +3.45%  '

# Row 29
$ws.Range("D29").Value = '18.85'
$ws.Range("E29").Value = '  +2.65%  '

# Row 30
$ws.Range("E30").Value = '  +1.12%  '

# Row 31
$ws.Range("D31").Value = '4.54'
$ws.Range("E31").Value = '  +6.53%  '

# Row 32
$ws.Range("D32").Value = '0.0610'
$ws.Range("E32").Value = '  +4.23%  '

# Row 33
$ws.Range("D33").Value = '1.96'
$ws.Range("E33").Value = '  +4.91%  '

# Row 34
$ws.Range("D34").Value = '4.32'
$ws.Range("E34").Value = '  +3.58%  '

# Row 35
$ws.Range("E35").Value = '  +0.04%  '

# Row 36
$ws.Range("D36").Value = '0.0839'
$ws.Range("E36").Value = '  +20.97%  '

# Row 37
$ws.Range("D37").Value = '1.48'
$ws.Range("E37").Value = '  -15.26%  '

# Row 38
$ws.Range("D38").Value = '0.861'
$ws.Range("E38").Value = '  +1.66%  '

# Row 39
$ws.Range("D39").Value = '2.01'
$ws.Range("E39").Value = '  +2.31%  '

# Row 40
$ws.Range("D40").Value = '103.99'
$ws.Range("E40").Value = '  +7.24%  '

# Row 41
$ws.Range("E41").Value = '  +3.92%  '

# Row 42
$ws.Range("D42").Value = '17.11'
$ws.Range("E42").Value = '  -0.90%  '

# Row 43
$ws.Range("D43").Value = '15.05'
$ws.Range("E43").Value = '  +21.92%  '

# Row 44
$ws.Range("D44").Value = '1.11'
$ws.Range("E44").Value = '  +2.63%  '

# Row 45
$ws.Range("D45").Value = '1.341.97'
$ws.Range("E45").Value = '  +3.30%  '

# Row 46
$ws.Range("D46").Value = '2.39'
$ws.Range("E46").Value = '  +2.60%  '

# Row 47
$ws.Range("D47").Value = '0.0808'
$ws.Range("E47").Value = '  +1.47%  '

# Row 48
$ws.Range("B48").Value = 'HuobiToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D48").Value = '2.41'
$ws.Range("E48").Value = '  +0.17%  '

# Row 49
$ws.Range("B49").Value = 'MXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D49").Value = '2.78'
$ws.Range("E49").Value = '  +1.68%  '

# Row 50
$ws.Range("D50").Value = '6.43'
$ws.Range("E50").Value = '  +2.77%  '

# Row 51
$ws.Range("D51").Value = '2.101.74'
$ws.Range("E51").Value = '  +1.33%  '
